$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4343.6797
$ws.Range("I132").Value = 2452.677
$ws.Range("K132").Value = 7358.031000000001
$ws.Range("M132").Value = -4828.031000000001
$ws.Range("H133").Value = 78800
$ws.Range("J133").Value = 77600
$ws.Range("L133").Value = 77600
$ws.Range("N133").Value = -87720
$ws.Range("H135").Value = 1770.4117
$ws.Range("I135").Value = 1000.9091
$ws.Range("K135").Value = 9008.1819
$ws.Range("M135").Value = -6473.1819
$ws.Range("H136").Value = 100071.336
$ws.Range("J136").Value = 100071.336
$ws.Range("L136").Value = 100071.336
$ws.Range("N136").Value = -110271.336
$ws.Range("H138").Value = 13159747
$ws.Range("I138").Value = 21740514
$ws.Range("K138").Value = 65221542
$ws.Range("M138").Value = -65216402

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 27234.666
$ws.Range("I33").Value = 19022.4
$ws.Range("K33").Value = 19022.4
$ws.Range("M33").Value = -18693.4
$ws.Range("H61").Value = 15771.143
$ws.Range("I61").Value = 15768.923
$ws.Range("K61").Value = 15768.923
$ws.Range("M61").Value = -15556.923
$ws.Range("H96").Value = 28554.545
$ws.Range("J96").Value = 28554.545
$ws.Range("L96").Value = 28554.545
$ws.Range("N96").Value = -34046.545
$ws.Range("H132").Value = 1123.22
$ws.Range("I132").Value = 1079.8334
$ws.Range("J132").Value = 1351
$ws.Range("K132").Value = 3239.5002
$ws.Range("L132").Value = 4053
$ws.Range("M132").Value = -709.5001999999999
$ws.Range("N132").Value = -9113
$ws.Range("H136").Value = 15771.143
$ws.Range("I136").Value = 15768.923
$ws.Range("K136").Value = 47306.769
$ws.Range("M136").Value = -44756.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1717.3383
$ws.Range("I20").Value = 1596.091
$ws.Range("J20").Value = 2230.3076
$ws.Range("K20").Value = 1596.091
$ws.Range("L20").Value = 2230.3076
$ws.Range("M20").Value = -1349.091
$ws.Range("N20").Value = -2724.3076
$ws.Range("H70").Value = 170964.5
$ws.Range("J70").Value = 170964.5
$ws.Range("L70").Value = 170964.5
$ws.Range("N70").Value = -171550.5
$ws.Range("H73").Value = 170964.5
$ws.Range("J73").Value = 170964.5
$ws.Range("L73").Value = 170964.5
$ws.Range("N73").Value = -172992.5
$ws.Range("H94").Value = 728.05
$ws.Range("I94").Value = 527.2941
$ws.Range("J94").Value = 1865.6666
$ws.Range("K94").Value = 527.2941
$ws.Range("L94").Value = 1865.6666
$ws.Range("M94").Value = -76.29409999999996
$ws.Range("N94").Value = -2767.6666
$ws.Range("H134").Value = 2596.0178
$ws.Range("I134").Value = 2308.8958
$ws.Range("J134").Value = 4318.75
$ws.Range("K134").Value = 6926.687399999999
$ws.Range("L134").Value = 12956.25
$ws.Range("M134").Value = -4391.687399999999
$ws.Range("N134").Value = -18026.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 3578.3845
$ws.Range("I3").Value = 4879.5
$ws.Range("J3").Value = 1496.6
$ws.Range("K3").Value = 4879.5
$ws.Range("L3").Value = 1496.6
$ws.Range("M3").Value = -4766.5
$ws.Range("N3").Value = -1722.6
$ws.Range("H31").Value = 2627.3076
$ws.Range("I31").Value = 1657.3334
$ws.Range("K31").Value = 1657.3334
$ws.Range("M31").Value = -1362.3334
$ws.Range("H33").Value = 530
$ws.Range("I33").Value = 530
$ws.Range("K33").Value = 530
$ws.Range("M33").Value = -151
$ws.Range("H34").Value = 2627.3076
$ws.Range("I34").Value = 1657.3334
$ws.Range("K34").Value = 1657.3334
$ws.Range("M34").Value = -1455.3334
$ws.Range("H58").Value = 1545.9474
$ws.Range("I58").Value = 1433.7646
$ws.Range("K58").Value = 1433.7646
$ws.Range("M58").Value = -1230.7646
$ws.Range("H94").Value = 2549.842
$ws.Range("I94").Value = 1408.375
$ws.Range("K94").Value = 1408.375
$ws.Range("M94").Value = -957.375
$ws.Range("H107").Value = 656.61536
$ws.Range("I107").Value = 450.9524
$ws.Range("J107").Value = 1520.4
$ws.Range("K107").Value = 450.9524
$ws.Range("L107").Value = 1520.4
$ws.Range("M107").Value = 1469.0476
$ws.Range("N107").Value = -5360.4
$ws.Range("H132").Value = 4700.069
$ws.Range("I132").Value = 2048.1226
$ws.Range("J132").Value = 19138.445
$ws.Range("K132").Value = 6144.3678
$ws.Range("L132").Value = 57415.335
$ws.Range("M132").Value = -3614.3678
$ws.Range("N132").Value = -62475.335
$ws.Range("H136").Value = 1545.9474
$ws.Range("I136").Value = 1433.7646
$ws.Range("K136").Value = 4301.293799999999
$ws.Range("M136").Value = -1751.293799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 153.26315
$ws.Range("J12").Value = 155.5
$ws.Range("L12").Value = 466.5
$ws.Range("N12").Value = -812.5
$ws.Range("H75").Value = 9999
$ws.Range("J75").Value = 9999
$ws.Range("L75").Value = 29997
$ws.Range("N75").Value = -31993
$ws.Range("H78").Value = 9999
$ws.Range("J78").Value = 9999
$ws.Range("L78").Value = 89991
$ws.Range("N78").Value = -99975
$ws.Range("H97").Value = 280.76923
$ws.Range("J97").Value = 348.55554
$ws.Range("L97").Value = 1045.66662
$ws.Range("N97").Value = -2037.66662
$ws.Range("H109").Value = 209.25
$ws.Range("I109").Value = 209.25
$ws.Range("K109").Value = 627.75
$ws.Range("M109").Value = 412.25
$ws.Range("H121").Value = 1626
$ws.Range("I121").Value = 1247.125
$ws.Range("K121").Value = 3741.375
$ws.Range("M121").Value = -2431.375
$ws.Range("H138").Value = 2180.3333
$ws.Range("I138").Value = 2180.3333
$ws.Range("K138").Value = 6540.999899999999
$ws.Range("M138").Value = -1400.999899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2574
$ws.Range("I132").Value = 2323.6416
$ws.Range("K132").Value = 6970.9248
$ws.Range("M132").Value = -4440.9248

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H95").Value = 41833.332
$ws.Range("J95").Value = 41833.332
$ws.Range("L95").Value = 41833.332
$ws.Range("N95").Value = -47325.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 94213.82000000001
$ws.Range("I96").Value = 103215.3
$ws.Range("J96").Value = 4199
$ws.Range("K96").Value = 103215.3
$ws.Range("L96").Value = 4199
$ws.Range("M96").Value = -101842.3
$ws.Range("N96").Value = -6945
$ws.Range("H126").Value = 11478.286
$ws.Range("I126").Value = 6369
$ws.Range("J126").Value = 24251.5
$ws.Range("K126").Value = 19107
$ws.Range("L126").Value = 72754.5
$ws.Range("M126").Value = -16637
$ws.Range("N126").Value = -77694.5
$ws.Range("H132").Value = 2136
$ws.Range("I132").Value = 1921.4314
$ws.Range("J132").Value = 3699.2856
$ws.Range("K132").Value = 5764.2942
$ws.Range("L132").Value = 11097.8568
$ws.Range("M132").Value = -3234.2942
$ws.Range("N132").Value = -16157.8568
$ws.Range("H135").Value = 239818.12
$ws.Range("J135").Value = 239818.12
$ws.Range("L135").Value = 239818.12
$ws.Range("N135").Value = -249958.12
$ws.Range("H136").Value = 1537.2554
$ws.Range("I136").Value = 1527.1957
$ws.Range("K136").Value = 4581.5871
$ws.Range("M136").Value = -2031.5871
